# Applies the crypto price/volume update described by the commit diff.
# A leading apostrophe forces Excel to store the assignment as literal text,
# matching the original inline-string cells (General format, type "s")
# instead of letting COM auto-coerce numeric-looking text into numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.195.51"
$ws.Range("E2").Value = "'  -0.82%  "

$ws.Range("D3").Value = "'1.880.14"
$ws.Range("E3").Value = "'  -1.67%  "

$ws.Range("E4").Value = "'  +0.06%  "

$ws.Range("D5").Value = "'237.08"
$ws.Range("E5").Value = "'  -0.79%  "

$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "'  +0.13%  "

$ws.Range("E7").Value = "'  -2.25%  "

$ws.Range("D8").Value = "'0.2795"
$ws.Range("E8").Value = "'  -2.24%  "

$ws.Range("D9").Value = "'0.06546"
$ws.Range("E9").Value = "'  -2.26%  "

$ws.Range("D10").Value = "'19.26"
$ws.Range("E10").Value = "'  +2.31%  "

$ws.Range("B11").Value = "'TRON"
$ws.Range("C11").Value = "'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").Value = "'0.07740"
$ws.Range("E11").Value = "'  +0.17%  "

$ws.Range("B12").Value = "'Litecoin"
$ws.Range("C12").Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D12").Value = "'98.09"
$ws.Range("E12").Value = "'  -4.56%  "

$ws.Range("D13").Value = "'1.896.93"
$ws.Range("E13").Value = "'  -0.78%  "

$ws.Range("D14").Value = "'5.108"
$ws.Range("E14").Value = "'  -1.85%  "

$ws.Range("D15").Value = "'0.6601"
$ws.Range("E15").Value = "'  -2.02%  "

$ws.Range("D16").Value = "'283.90"
$ws.Range("E16").Value = "'  +9.76%  "

$ws.Range("D17").Value = "'30.162.87"
$ws.Range("E17").Value = "'  -0.91%  "

$ws.Range("D18").Value = "'1.000"
$ws.Range("E18").Value = "'  -0.01%  "

$ws.Range("D19").Value = "'2.144.06"
$ws.Range("E19").Value = "'  -0.55%  "

$ws.Range("D20").Value = "'12.41"
$ws.Range("E20").Value = "'  -2.50%  "

$ws.Range("D21").Value = "'0.000007254"
$ws.Range("E21").Value = "'  -3.41%  "

$ws.Range("D22").Value = "'5.300"
$ws.Range("E22").Value = "'  -2.02%  "

$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "'  +0.06%  "

$ws.Range("D24").Value = "'6.168"
$ws.Range("E24").Value = "'  -2.24%  "

$ws.Range("D25").Value = "'9.225"
$ws.Range("E25").Value = "'  -2.57%  "

$ws.Range("E26").Value = "'  +0.37%  "

$ws.Range("D27").Value = "'18.89"
$ws.Range("E27").Value = "'  -0.63%  "

$ws.Range("D28").Value = "'1.993"
$ws.Range("E28").Value = "'  -3.42%  "

$ws.Range("D29").Value = "'1.385"
$ws.Range("E29").Value = "'  +0.62%  "

$ws.Range("D30").Value = "'0.09792"
$ws.Range("E30").Value = "'  -3.07%  "

$ws.Range("D31").Value = "'4.439"
$ws.Range("E31").Value = "'  -4.35%  "

$ws.Range("D32").Value = "'1.487"
$ws.Range("E32").Value = "'  -1.69%  "

$ws.Range("D33").Value = "'4.168"
$ws.Range("E33").Value = "'  -2.35%  "

$ws.Range("D34").Value = "'0.04652"
$ws.Range("E34").Value = "'  -2.93%  "

$ws.Range("D35").Value = "'0.7044"
$ws.Range("E35").Value = "'  -3.92%  "

$ws.Range("D36").Value = "'1.088"
$ws.Range("E36").Value = "'  -2.35%  "

$ws.Range("E37").Value = "'  +0.17%  "

$ws.Range("D38").Value = "'2.715"
$ws.Range("E38").Value = "'  +0.14%  "

$ws.Range("D39").Value = "'0.01853"
$ws.Range("E39").Value = "'  -4.14%  "

$ws.Range("D40").Value = "'6.726"
$ws.Range("E40").Value = "'  +7.47%  "

$ws.Range("D41").Value = "'2.523"
$ws.Range("E41").Value = "'  -2.91%  "

$ws.Range("D42").Value = "'72.30"
$ws.Range("E42").Value = "'  -3.11%  "

$ws.Range("D43").Value = "'0.8685"
$ws.Range("E43").Value = "'  +0.26%  "

$ws.Range("D44").Value = "'1.942"
$ws.Range("E44").Value = "'  -2.98%  "

$ws.Range("D45").Value = "'1.002"
$ws.Range("E45").Value = "'  +0.13%  "

$ws.Range("E46").Value = "'  -2.65%  "

$ws.Range("D47").Value = "'0.4150"
$ws.Range("E47").Value = "'  -2.60%  "

$ws.Range("D48").Value = "'997.67"
$ws.Range("E48").Value = "'  -7.12%  "

$ws.Range("D49").Value = "'7.180"
$ws.Range("E49").Value = "'  -4.17%  "

$ws.Range("D50").Value = "'9.144"
$ws.Range("E50").Value = "'  +3.57%  "

$ws.Range("D51").Value = "'0.1161"
$ws.Range("E51").Value = "'  -3.30%  "
